$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 294800
$ws.Range("E8").Value = 249200
$ws.Range("F8").Value = 213000
$ws.Range("G8").Value = 220700
$ws.Range("H8").Value = 335800
$ws.Range("I8").Value = 286400
$ws.Range("J8").Value = 2134400
$ws.Range("D9").Value = 250500
$ws.Range("E9").Value = 223700
$ws.Range("F9").Value = 197200
$ws.Range("G9").Value = 399800
$ws.Range("H9").Value = 596500
$ws.Range("I9").Value = 467900
$ws.Range("J9").Value = 1707800
$ws.Range("D10").Value = 44200
$ws.Range("E10").Value = 25500
$ws.Range("F10").Value = 15800
$ws.Range("G10").Value = -179100
$ws.Range("H10").Value = -260700
$ws.Range("I10").Value = -181500
$ws.Range("J10").Value = 426500
$ws.Range("J14").Value = 329100
$ws.Range("J15").Value = 8800
$ws.Range("D17").Value = 269400
$ws.Range("E17").Value = 239100
$ws.Range("F17").Value = 211500
$ws.Range("G17").Value = 215900
$ws.Range("H17").Value = 300500
$ws.Range("I17").Value = 258800
$ws.Range("J17").Value = 2274500
$ws.Range("D18").Value = 25300
$ws.Range("E18").Value = 10200
$ws.Range("F18").Value = 1500
$ws.Range("G18").Value = 4800
$ws.Range("H18").Value = 35300
$ws.Range("I18").Value = 27600
$ws.Range("J18").Value = -140200
$ws.Range("D20").Value = -5600
$ws.Range("F20").Value = 131500
$ws.Range("G20").Value = 7000
$ws.Range("H20").Value = -4100
$ws.Range("I20").Value = -6600
$ws.Range("J20").Value = -142900
$ws.Range("D21").Value = 22500
$ws.Range("E21").Value = 13100
$ws.Range("F21").Value = 136100
$ws.Range("G21").Value = 14100
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 98800
$ws.Range("G22").Value = 22900
$ws.Range("H22").Value = 24600
$ws.Range("I22").Value = 31700
$ws.Range("D23").Value = 15000
$ws.Range("E23").Value = 5700
$ws.Range("F23").Value = 127300
$ws.Range("G23").Value = -11000
$ws.Range("H23").Value = 6600
$ws.Range("I23").Value = -10800
$ws.Range("J23").Value = -283100
$ws.Range("D24").Value = 3500
$ws.Range("G24").Value = 11800
$ws.Range("H24").Value = 23200
$ws.Range("I24").Value = 9400
$ws.Range("J24").Value = 51900
$ws.Range("D26").Value = 11500
$ws.Range("E26").Value = 5700
$ws.Range("F26").Value = 126100
$ws.Range("G26").Value = -22800
$ws.Range("H26").Value = -16600
$ws.Range("I26").Value = -20100
$ws.Range("J26").Value = -334900
$ws.Range("D27").Value = 11500
$ws.Range("E27").Value = 5700
$ws.Range("F27").Value = 126100
$ws.Range("G27").Value = -22800
$ws.Range("H27").Value = 30300
$ws.Range("I27").Value = 536100
$ws.Range("J27").Value = -149000
$ws.Range("H29").Value = -70400
$ws.Range("I29").Value = -770900
$ws.Range("D32").Value = 5600
$ws.Range("F32").Value = -131500
$ws.Range("G32").Value = -7000
$ws.Range("H32").Value = 4100
$ws.Range("I32").Value = 6600
$ws.Range("J32").Value = 142900
$ws.Range("D33").Value = 11500
$ws.Range("E33").Value = 5700
$ws.Range("F33").Value = 126100
$ws.Range("G33").Value = -22800
$ws.Range("H33").Value = -40000
$ws.Range("I33").Value = -234800
$ws.Range("J33").Value = -149000
$ws.Range("D35").Value = 11500
$ws.Range("E35").Value = 5700
$ws.Range("F35").Value = 126100
$ws.Range("G35").Value = -22800
$ws.Range("H35").Value = -40000
$ws.Range("I35").Value = -234800
$ws.Range("J35").Value = -149000
$ws.Range("D41").Value = 23100
$ws.Range("E41").Value = 9400
$ws.Range("F41").Value = 6700
$ws.Range("J41").Value = 294100
$ws.Range("G42").Value = 92700
$ws.Range("J42").Value = 37500
$ws.Range("D43").Value = 27700
$ws.Range("E43").Value = 33100
$ws.Range("F43").Value = 33600
$ws.Range("G43").Value = 21800
$ws.Range("H43").Value = 33100
$ws.Range("J43").Value = 477600
$ws.Range("D44").Value = 26300
$ws.Range("E44").Value = 16700
$ws.Range("F44").Value = 36300
$ws.Range("G44").Value = 18500
$ws.Range("H44").Value = 36400
$ws.Range("J44").Value = 59300
$ws.Range("D45").Value = 22300
$ws.Range("E45").Value = 11100
$ws.Range("F45").Value = 9400
$ws.Range("H45").Value = 35600
$ws.Range("J45").Value = 5200
$ws.Range("D46").Value = 99400
$ws.Range("E46").Value = 70300
$ws.Range("F46").Value = 86800
$ws.Range("G46").Value = 138200
$ws.Range("H46").Value = 106500
$ws.Range("J46").Value = 873800
$ws.Range("F47").Value = 1900
$ws.Range("H47").Value = 105400
$ws.Range("J47").Value = 236400
$ws.Range("D48").Value = 4800
$ws.Range("E48").Value = 4900
$ws.Range("F48").Value = 5000
$ws.Range("H48").Value = 4400
$ws.Range("J48").Value = 566100
$ws.Range("D49").Value = 6400
$ws.Range("E49").Value = 8600
$ws.Range("F49").Value = 10700
$ws.Range("G49").Value = 13200
$ws.Range("H49").Value = 15200
$ws.Range("J49").Value = 2394500
$ws.Range("H52").Value = 13000
$ws.Range("J52").Value = 77300
$ws.Range("D54").Value = 110900
$ws.Range("E54").Value = 85700
$ws.Range("F54").Value = 104400
$ws.Range("G54").Value = 157300
$ws.Range("H54").Value = 244400
$ws.Range("J54").Value = 4148200
$ws.Range("D57").Value = 5500
$ws.Range("E57").Value = 35300
$ws.Range("F57").Value = 55200
$ws.Range("G57").Value = 9900
$ws.Range("H57").Value = 73900
$ws.Range("J57").Value = 315600
$ws.Range("D58").Value = 22600
$ws.Range("G58").Value = 263800
$ws.Range("H58").Value = 86900
$ws.Range("J58").Value = 339900
$ws.Range("D59").Value = 4900
$ws.Range("E59").Value = 4300
$ws.Range("G59").Value = 8000
$ws.Range("H59").Value = 13500
$ws.Range("J59").Value = 102900
$ws.Range("D60").Value = 33100
$ws.Range("E60").Value = 42400
$ws.Range("F60").Value = 64200
$ws.Range("G60").Value = 281700
$ws.Range("H60").Value = 174400
$ws.Range("J60").Value = 758400
$ws.Range("D61").Value = 33500
$ws.Range("E61").Value = 33100
$ws.Range("F61").Value = 33300
$ws.Range("H61").Value = 177100
$ws.Range("J61").Value = 1945400
$ws.Range("F62").Value = 8500
$ws.Range("G62").Value = 11600
$ws.Range("H62").Value = 13200
$ws.Range("J62").Value = 170800
$ws.Range("D66").Value = 67500
$ws.Range("E66").Value = 78500
$ws.Range("F66").Value = 106400
$ws.Range("G66").Value = 293600
$ws.Range("H66").Value = 365000
$ws.Range("J66").Value = 3987900
$ws.Range("D72").Value = 213400
$ws.Range("E72").Value = 198400
$ws.Range("F72").Value = 192500
$ws.Range("G72").Value = 65900
$ws.Range("H72").Value = 85000
$ws.Range("J72").Value = 361200
$ws.Range("D76").Value = 43400
$ws.Range("E76").Value = 7200
$ws.Range("G76").Value = -136300
$ws.Range("H76").Value = -120600
$ws.Range("J76").Value = 160300
$ws.Range("D81").Value = 11500
$ws.Range("E81").Value = 5700
$ws.Range("F81").Value = 126100
$ws.Range("G81").Value = -22800
$ws.Range("H81").Value = -40000
$ws.Range("I81").Value = -234800
$ws.Range("J81").Value = -149000
$ws.Range("D83").Value = 2700
$ws.Range("J83").Value = 382700
$ws.Range("D89").Value = -31800
$ws.Range("E89").Value = 8000
$ws.Range("F89").Value = 13400
$ws.Range("G89").Value = -35300
$ws.Range("H89").Value = 80600
$ws.Range("J89").Value = 484200
$ws.Range("E91").Value = -400
$ws.Range("J91").Value = -96300
$ws.Range("F94").Value = -6100
$ws.Range("G94").Value = 54600
$ws.Range("H94").Value = 58200
$ws.Range("J94").Value = -312300
$ws.Range("J96").Value = -27600
$ws.Range("D100").Value = 40400
$ws.Range("E100").Value = -3900
$ws.Range("G100").Value = -19300
$ws.Range("H100").Value = -150400
$ws.Range("J100").Value = -27300
$ws.Range("J101").Value = -1900
$ws.Range("D102").Value = 13600
$ws.Range("E102").Value = 2700
$ws.Range("H102").Value = -12700
$ws.Range("J102").Value = 142600
